$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Activity Log")

# Row 28: last4=6977, date=2020-04-03 (43924), start=23:15 (0.96875), end=23:45 (0.98958333333333337)
$ws.Range("B28").Value = 6977
$ws.Range("C28").Value = 43924
$ws.Range("D28").Value = 0.96875
$ws.Range("E28").Value = 0.98958333333333337
$ws.Range("G28").Value = "Generated the functional waveforms of the LogicUnit.vhd. Exported them into Documentation as per instructions"

# Row 29: last4=6977, date=2020-04-03 (43924), start=23:45 (0.98958333333333337), end=23:55 (0.99652777777777779)
$ws.Range("B29").Value = 6977
$ws.Range("C29").Value = 43924
$ws.Range("D29").Value = 0.98958333333333337
$ws.Range("E29").Value = 0.99652777777777779
$ws.Range("G29").Value = "Systhesied circuits but not satisfied with the diagrams due to it being very cluttered. Will grab the images later when a revised version is done"

# Row 30: last4=6977, date=2020-04-03 (43924), start=23:55 (0.99652777777777779), end=00:05 (3.472222222222222E-3)
$ws.Range("B30").Value = 6977
$ws.Range("C30").Value = 43924
$ws.Range("D30").Value = 0.99652777777777779
$ws.Range("E30").Value = 0.003472222222222222
$ws.Range("G30").Value = "Set up files and environment to obtain timing simulations from ModelSim"

# Row 31: last4=6977, date=2020-04-03 (43924), start=00:05 (3.472222222222222E-3), end=00:07 (4.8611111111111112E-3)
$ws.Range("B31").Value = 6977
$ws.Range("C31").Value = 43924
$ws.Range("D31").Value = 0.003472222222222222
$ws.Range("E31").Value = 0.004861111111111111
$ws.Range("G31").Value = "Updated .gitignore to ignore temporary files that is unneeded"

# Update the active selection to match the new cursor position
$ws.Range("A31").Select()
